$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (rows 2..64 hold player records)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if (-not $lastRow -or $lastRow -lt 64) { $lastRow = 64 }

# New header cells, matching the style of the existing header row
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every player row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 52
    $ws.Cells.Item($r, 31).Value = 110
    $ws.Cells.Item($r, 32).Value = 0
}
